$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Z-czni")

# New column header in F1 ("Odsetek")
$ws.Range("F1").Value = "Odsetek"

# Row 2 holds the "anchor" formulas
$ws.Range("E2").Formula = "=C2+D2"
$ws.Range("F2").Formula = "=C2/E2*100"

# Rows 3:41 filled in one go so the engine stores them as shared formulas,
# matching the fill-down Excel would perform
$ws.Range("E3:E41").Formula = "=C3+D3"
$ws.Range("F3:F41").Formula = "=C3/E3*100"

# Reflect the new column in the sheet view / selection
$ws.Activate()
$ws.Range("F2:F41").Select()
